$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The tracker sheet keeps one column per date snapshot, newest first (after
# the fixed "Company"/A column). This revision adds two newer snapshots
# ("Jun_17" and "Jun_15") ahead of the existing "Jun_13" column, so insert
# two new columns at C - this pushes the old "Jun_10" column (C) to E and
# leaves the old "Jun_13" column letter (B) in place structurally.
$ws.Range("C1").EntireColumn.Insert()
$ws.Range("C1").EntireColumn.Insert()

# Give the two newly inserted columns the same on-screen width as their
# neighbours (42.5 "characters" in the saved XML == 41.666... in the
# Application.ColumnWidth units COM exposes).
$ws.Range("C1").ColumnWidth = 41.666666666666664
$ws.Range("D1").ColumnWidth = 41.666666666666664

# Row 1 holds the date-snapshot headers. The newest dates go in B1/C1, and
# the previous header ("Jun_13") shifts from B1 into the newly freed D1.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = "Jun_13"

# The new snapshot columns have no rating data yet, so every data row gets
# the same "UN" (unrated) placeholder used elsewhere on the sheet.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}
